$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""66.813.95"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = "=""  +4.30%  """
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = "=""3.232.26"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = "=""  +1.94%  """
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("E4").Formula = "=""  +0.00%  """
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("D5").Formula = "=""575.57"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = "=""  +3.44%  """
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""175.31"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = "=""  +3.09%  """
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = "=""  -0.06%  """
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("E8").Formula = "=""  +0.44%  """
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("D9").Formula = "=""3.228.62"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = "=""  +1.94%  """
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("E10").Formula = "=""  +4.88%  """
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""6.66"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = "=""  +1.44%  """
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("E12").Formula = "=""  +2.65%  """
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=""3.790.90"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = "=""  +1.97%  """
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("E14").Formula = "=""  +2.08%  """
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("D15").Formula = "=""27.63"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = "=""  +1.34%  """
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=""66.796.10"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = "=""  +3.97%  """
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("E17").Formula = "=""  +4.35%  """
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D18").Formula = "=""3.226.25"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = "=""  +1.60%  """
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("D19").Formula = "=""5.77"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = "=""  +2.82%  """
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=""13.18"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = "=""  +2.12%  """
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("D21").Formula = "=""365.26"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = "=""  +3.98%  """
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("E22").Formula = "=""  +3.70%  """
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("E23").Formula = "=""  +0.30%  """
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=""69.70"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = "=""  +1.73%  """
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""0.0000118"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = "=""  +0.80%  """
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("E26").Formula = "=""  +0.54%  """
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""3.356.33"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = "=""  +0.91%  """
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""9.82"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = "=""  +4.92%  """
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""0.178"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = "=""  +2.28%  """
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("E30").Formula = "=""  +0.17%  """
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("E31").Formula = "=""  +3.56%  """
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""5.56"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = "=""  -0.41%  """
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("B33").Formula = "=""USDe"""
$ws.Range("B33").Copy()
$ws.Range("B33").PasteSpecial(-4163)
$ws.Range("C33").Formula = "=""https://coinranking.com/coin/exbfr2U-0+usde-usde"""
$ws.Range("C33").Copy()
$ws.Range("C33").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""0.998"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = "=""  -0.18%  """
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("B34").Formula = "=""EthereumClassic"""
$ws.Range("B34").Copy()
$ws.Range("B34").PasteSpecial(-4163)
$ws.Range("C34").Formula = "=""https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"""
$ws.Range("C34").Copy()
$ws.Range("C34").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=""22.31"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = "=""  +1.26%  """
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("B35").Formula = "=""Fetch.AI"""
$ws.Range("B35").Copy()
$ws.Range("B35").PasteSpecial(-4163)
$ws.Range("C35").Formula = "=""https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"""
$ws.Range("C35").Copy()
$ws.Range("C35").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=""1.22"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = "=""  +3.47%  """
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("B36").Formula = "=""Aptos"""
$ws.Range("B36").Copy()
$ws.Range("B36").PasteSpecial(-4163)
$ws.Range("C36").Formula = "=""https://coinranking.com/coin/HGYj5JCv5+aptos-apt"""
$ws.Range("C36").Copy()
$ws.Range("C36").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=""6.72"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = "=""  +2.34%  """
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("D37").Formula = "=""168.20"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = "=""  +6.46%  """
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("E38").Formula = "=""  +4.51%  """
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("E39").Formula = "=""  +4.15%  """
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""1.84"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = "=""  +10.17%  """
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=""26.62"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = "=""  +1.94%  """
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""2.60"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = "=""  +4.63%  """
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=""6.36"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = "=""  +6.07%  """
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=""2.680.99"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = "=""  +1.51%  """
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("E45").Formula = "=""  +3.12%  """
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("E46").Formula = "=""  +4.48%  """
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=""0.0670"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = "=""  +3.09%  """
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=""24.39"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = "=""  +5.24%  """
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=""329.97"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = "=""  +3.52%  """
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("E50").Formula = "=""  +3.18%  """
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("E51").Formula = "=""  +1.16%  """
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
